$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.961.95"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.818.12"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.93"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.50"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.50"
$ws.Range("E10").Value = "  -6.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.88"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.74"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.258.99"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.812.36"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.932"
$ws.Range("E17").Value = "  +4.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.803.29"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").Value = "  +5.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.14"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.42"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0992"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.56"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.61"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.92"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.31"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0491"
$ws.Range("E30").Value = "  +19.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.144"
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.59"
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.61"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.53"
$ws.Range("E35").Value = "  +10.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0850"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -4.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.37"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.118"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.93"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.07"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("E44").Value = "  -7.84%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.072.90"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.85"
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.969"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.12"
$ws.Range("E51").Value = "  +1.58%  "
